# "improve master plan -> AI OT Manager Testing III (20241113)"
#
# The master plan test data was regenerated: rows 51-66 (the last 16 cases,
# which had rolled forward to Jan/Feb 2025 booking & op dates) were dropped,
# and the remaining cases (rows 2-50) had their BOOKING DATE / OPERATION DATE
# (cols A/B) and OT LIST NAME sequence number (col M) refreshed to new values.
# The view was also left scrolled back to the top with a new zoom level.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New BOOKING DATE (A) / OPERATION DATE (B) / OT LIST NAME # (M) values ---
# Columns: row, A (booking date serial), B (operation date serial), M (seq #)
$rowData = @(
    @(2,45546,44929,12),
    @(3,44880,44930,6),
    @(4,44861,44929,7),
    @(5,44917,44930,1),
    @(6,44924,44929,3),
    @(7,44916,44929,8),
    @(8,44923,44929,7),
    @(9,44917,44929,7),
    @(10,44923,44929,6),
    @(11,44929,44929,5),
    @(12,44837,44930,7),
    @(13,44893,44930,5),
    @(14,44895,44930,4),
    @(15,44896,44930,3),
    @(16,44861,44929,3),
    @(17,44895,44929,2),
    @(18,44929,44929,5),
    @(19,44910,44930,1),
    @(20,44847,44929,5),
    @(21,44929,44929,13),
    @(22,44917,44929,7),
    @(23,44922,44929,8),
    @(24,44926,44929,10),
    @(25,44929,44929,7),
    @(26,44929,44929,4),
    @(27,44929,44929,6),
    @(28,44929,44929,6),
    @(29,44925,44929,13),
    @(30,44866,44930,1),
    @(31,44925,44929,9),
    @(32,44914,44929,5),
    @(33,44918,44929,3),
    @(34,44918,44929,2),
    @(35,44925,44929,11),
    @(36,44922,44929,4),
    @(37,44812,44929,2),
    @(38,44929,44929,9),
    @(39,44928,44928,5),
    @(40,44886,44929,9),
    @(41,44924,44929,7),
    @(42,44866,44930,8),
    @(43,44889,44929,1),
    @(44,44929,44929,7),
    @(45,44923,44929,6),
    @(46,44929,44929,11),
    @(47,44923,44929,1),
    @(48,44927,44929,10),
    @(49,44924,44929,6),
    @(50,44908,44930,2)
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $ws.Range("A$r").Value = $entry[1]
    $ws.Range("B$r").Value = $entry[2]
    $ws.Range("M$r").Value = $entry[3]
}

# --- Drop the trailing cases (rows 51-66); this also prunes the now-unused
#     "N290988" shared string and shrinks the sheet dimension to A1:Q50 ---
$ws.Rows("51:66").Delete()

# --- Column K (11) is a touch narrower now ---
$ws.Columns(11).ColumnWidth = 26.6

# --- View scrolled back to the top, zoomed to 83%, with M7 selected ---
$excel.ActiveWindow.Zoom = 83
$ws.Range("M7").Select()

Write-Output "masterplan refreshed: 49 data rows (was 65), view reset"
